{"js": "// The document contains three words that use the letter \"\u0451\" (e with\n// diaeresis/umlaut): \"\u0443\u0447\u0451\u0442\u043d\u0430\u044f\", \"\u0432\u0432\u0435\u0434\u0451\u043d\u043d\u044b\u0445\", \"\u0441\u043e\u0437\u0434\u0430\u0451\u0442\". The edit replaces\n// each of those occurrences with the plain letter \"\u0435\", turning the words\n// into \"\u0443\u0447\u0435\u0442\u043d\u0430\u044f\", \"\u0432\u0432\u0435\u0434\u0435\u043d\u043d\u044b\u0445\", \"\u0441\u043e\u0437\u0434\u0430\u0435\u0442\" respectively, while leaving the\n// rest of each sentence untouched.\nconst body = context.document.body;\n\n// Find every occurrence of \"\u0451\" in the document body and swap it for \"\u0435\".\nconst matches = body.search(\"\u0451\", { matchCase: true });\nmatches.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < matches.items.length; i++) {\n  matches.items[i].insertText(\"\u0435\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# The document contains three words spelled with the letter \"\u0451\" (e with\n# diaeresis): \"\u0443\u0447\u0451\u0442\u043d\u0430\u044f\", \"\u0432\u0432\u0435\u0434\u0451\u043d\u043d\u044b\u0445\", \"\u0441\u043e\u0437\u0434\u0430\u0451\u0442\". The edit swaps that\n# letter for the plain \"\u0435\" in each of them (-> \"\u0443\u0447\u0435\u0442\u043d\u0430\u044f\", \"\u0432\u0432\u0435\u0434\u0435\u043d\u043d\u044b\u0445\",\n# \"\u0441\u043e\u0437\u0434\u0430\u0435\u0442\"), leaving the rest of the surrounding sentences unchanged.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"\u0451\"\n$find.Replacement.Text = \"\u0435\"\n$find.Forward = $true\n$find.Wrap = 1            # wdFindContinue - keep searching the whole story\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n# wdReplaceAll = 2 -> replace every remaining match of \"\u0451\" with \"\u0435\"\n$find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $false, $false, $find.Forward, $find.Wrap, $find.Format, $find.Replacement.Text, 2)\n"}
